$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> B, C, D values (B uses $null when the cell should stay empty)
$data = @(
    @{Row=2;  B=54000; C=4944.292321154038; D=49055.70767884596},
    @{Row=3;  B=50801; C=4887.559847582665; D=45913.44015241734},
    @{Row=4;  B=48611; C=4825.545484664752; D=43785.45451533525},
    @{Row=5;  B=46724; C=4794.057447783162; D=41929.94255221684},
    @{Row=6;  B=47566; C=4794.010622970725; D=42771.98937702928},
    @{Row=7;  B=48485; C=4890.506649537315; D=43594.49335046268},
    @{Row=8;  B=47871; C=5267.449152182856; D=42603.55084781715},
    @{Row=9;  B=41983; C=6573.678702600464; D=35409.32129739954},
    @{Row=10; B=65400; C=8191.789500000001; D=57208.2105},
    @{Row=11; B=73839; C=14018.138;         D=59820.862},
    @{Row=12; B=77636; C=14986.9245;        D=62649.0755},
    @{Row=13; B=78566; C=14980.626;         D=63585.374},
    @{Row=14; B=79318; C=15169.25375;       D=64148.74625},
    @{Row=15; B=83095; C=15666.69225;       D=67428.30775000001},
    @{Row=16; B=51910; C=15813.73925;       D=36096.26075},
    @{Row=17; B=$null; C=16092.652;         D=0},
    @{Row=18; B=$null; C=16268.66525;       D=1768.269410000001},
    @{Row=19; B=$null; C=16049.11175;       D=0},
    @{Row=20; B=$null; C=14955.2725;        D=0},
    @{Row=21; B=$null; C=13776.544;         D=0},
    @{Row=22; B=$null; C=12164.7975;        D=0},
    @{Row=23; B=$null; C=9625.974749999999; D=0},
    @{Row=24; B=$null; C=6491.415499999999; D=0},
    @{Row=25; B=$null; C=5924.492750000001; D=0}
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.B -ne $null) {
        $ws.Cells.Item($r, 2).Value = $item.B
    }
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}
